$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update account number on row 3 (H3: 4672846545 -> 7166033480)
$ws.Range("H3").Value = 7166033480

# 2. Update the FechaInicio shared text value used by T3/T4/T5 (27/03/2023 -> 14/09/2020)
$ws.Range("T3").Value = "'14/09/2020"

# 3. Populate row 4 by copying the pattern from row 3, then adjusting the
#    sequence number and account number. Copy in two pieces (A:M and R:T)
#    so we don't materialise blank cells in the untouched N:Q gap.
$ws.Range("A3:M3").Copy($ws.Range("A4:M4"))
$ws.Range("R3:T3").Copy($ws.Range("R4:T4"))
$ws.Range("A4").Value = 3
$ws.Range("H4").Value = 7166033480

# 4. Populate row 5 the same way.
$ws.Range("A3:M3").Copy($ws.Range("A5:M5"))
$ws.Range("R3:T3").Copy($ws.Range("R5:T5"))
$ws.Range("A5").Value = 4
$ws.Range("H5").Value = 7166033480

# 5. Add the hyperlinks for the new C4 / C5 cells (same target URL as C3)
#    and restore the "Hipervínculo" cell style that Hyperlinks.Add doesn't
#    apply by itself.
$ws.Hyperlinks.Add($ws.Range("C4"), "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/pc/PolicyCenter.do")
$ws.Range("C4").Style = "Hipervínculo"
$ws.Hyperlinks.Add($ws.Range("C5"), "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/pc/PolicyCenter.do")
$ws.Range("C5").Style = "Hipervínculo"

# 6. Reset the view: scroll back to A1 (drop topLeftCell="F1") and move the
#    active selection to A3.
$ws.Range("A1").Select()
$ws.Range("A3").Select()
